$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- New tracklist data for Sheet1 (rows 2-15) ---
# Columns: B=Title, C=Composer, D=Performer, E=Time (fraction of day)
$titles = @(
    "Chandler",
    "Falling in Love",
    "On One",
    "Losing",
    "You",
    "Westside Gunn's Interlude",
    "Addicted",
    "Kaytra's Interlude",
    "Must Be Nice",
    "Hot Minute Interlude",
    "Think About You",
    "So So Sick",
    "That's on You [Japanese Remix]",
    "Overgrown"
)
$composers = @(
    "Joyce Wrice, Mack Keane, Preston Harris",
    "Daniel Church, Davion Farris, Joyce Wrice, Lucy Daye",
    "Daniel Church, Freddie Gibbs, Joyce Wrice",
    "Daniel Church, Joyce Wrice",
    "Daniel Church, Joyce Wrice",
    "Westside Gunn",
    "D'Mile, Daniel Church, James Poyser, Jonah Christian, Joyce Wrice",
    "Joyce Wrice, KATRANADA, Mack Keane",
    "D'Mile, Jonah Christian, Joyce Wrice, Masego, Shawn",
    "Devin Morrison, Joyce Wrice, Mndsgn",
    "Bradford Tidwell, Jason Kawu-Eugenio, Jonah Christian, Joyce Wrice, Varren Wade",
    "Daniel Church, Joyce Wrice",
    "Austin Brown, Joyce Wrice, Sol Was, UMI",
    "Joyce Wrice, Mack Keane"
)
$performers = @(
    "Joyce Wrice",
    "Joyce Wrice feat: Lucky Daye",
    "Joyce Wrice feat. Freddie Gibbs",
    "Joyce Wrice",
    "Joyce Wrice",
    "Joyce Wrice feat: Westside Gunn, ESTA.",
    "Joyce Wrice",
    "Joyce Wrice feat: KAYTRANADA",
    "Joyce Wrice feat: Masego",
    "Joyce Wrice feat: Mndsgn, Devin Morrison",
    "Joyce Wrice",
    "Joyce Wrice",
    "Joyce Wrice Feat: UMI",
    "Joyce Wrice"
)
$times = @(
    0.1277777777777778,
    0.12430555555555556,
    0.14097222222222222,
    0.13194444444444445,
    0.07569444444444444,
    0.08541666666666665,
    0.1361111111111111,
    0.05555555555555555,
    0.14930555555555555,
    0.05277777777777778,
    0.11666666666666665,
    0.16597222222222222,
    0.11597222222222221,
    0.13263888888888889
)

for ($i = 0; $i -lt $titles.Count; $i++) {
    $row = 2 + $i
    foreach ($ws in @($ws1, $ws3)) {
        $ws.Cells.Item($row, 1).Value2 = $i + 1
        $ws.Cells.Item($row, 2).Value2 = $titles[$i]
        $ws.Cells.Item($row, 3).Value2 = $composers[$i]
        $ws.Cells.Item($row, 4).Value2 = $performers[$i]
        $ws.Cells.Item($row, 5).Value2 = $times[$i]
    }
}

# --- Column widths on Sheet1 / Sheet3 ---
foreach ($ws in @($ws1, $ws3)) {
    $ws.Columns.Item(2).ColumnWidth = 29.77734375
    $ws.Columns.Item(3).ColumnWidth = 74.77734375
    $ws.Columns.Item(4).ColumnWidth = 38.88671875
    $ws.Columns.Item(5).ColumnWidth = 8.44140625
}

# --- Recalculate so Sheet2's report formulas pick up new values ---
$wb.Application.Calculate()

# --- Update selection on Sheet2 to reflect the new data extent ---
$ws2.Activate()
$ws2.Range("A3:K18").Select()

# --- Defined names: rename serpentwithfeet2 -> joycewrice1 and widen range ---
$wb.Names.Item("Sheet1!serpentwithfeet2").Name = "joycewrice1"
$wb.Names.Item("Sheet1!joycewrice1").RefersTo = "=Sheet1!`$A`$1:`$E`$15"
$wb.Names.Item("Sheet3!serpentwithfeet2").Name = "joycewrice1"
$wb.Names.Item("Sheet3!joycewrice1").RefersTo = "=Sheet3!`$A`$1:`$E`$15"
